# TOR-016 New Band Implementation - data refresh for Band.xlsx
# Applies:
#  - shared-string text fixes (trim trailing spaces, shorten a name, fix a typo)
#  - Bands sheet: eyeAwareness/shadowPts value updates + alignment/protection
#    formatting on the hopePts/shadowPts/shadowScars columns
#  - Allies sheet: populate the new "kinglyGift" column with NONE, turn on
#    word-wrap for the "quirksOrNotes" column body rows, shrink row 7's
#    custom height now that the typo-fixed note still wraps, and fix the
#    "Eiri'ego" -> "Eitri'ego" typo inside that note.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Bands sheet
# ---------------------------------------------------------------------
$bands = $wb.Worksheets.Item("Bands")

# Turn on alignment/protection formatting for the hopePts/shadowPts/shadowScars
# header + data cells (I:K), matching the rest of the header row.
$bands.Range("I1:K2").HorizontalAlignment = -4130

# Data updates for Veig's Band (row 2)
$bands.Range("F2").Value = 3
$bands.Range("J2").Value = 2

$bands.Range("H10").Select()

# ---------------------------------------------------------------------
# Allies sheet
# ---------------------------------------------------------------------
$allies = $wb.Worksheets.Item("Allies")

# Trim trailing spaces from header labels
$allies.Range("H1").Value = "gift"
$allies.Range("L1").Value = "quirksOrNotes"

# Shorten Galar's display name
$allies.Range("C3").Value = "Galar"

# Fix typo in Regin's quirks/notes text
$allies.Range("L7").Value = "Tarczę zdobi wojennymi trofeami:" + [char]10 + "- Pot. pod Aleją Królów" + [char]10 + "- Ekstrakcja inż. Eitri'ego" + [char]10 + "Walczy włócznią z tarczą"

# Populate the new "kinglyGift" column (J) for every ally with NONE
$allies.Range("J2").Value = "NONE"
$allies.Range("J3").Value = "NONE"
$allies.Range("J4").Value = "NONE"
$allies.Range("J5").Value = "NONE"
$allies.Range("J6").Value = "NONE"
$allies.Range("J7").Value = "NONE"
$allies.Range("J8").Value = "NONE"

# Turn on word-wrap for the quirksOrNotes body rows (header/row7 already wrap)
$allies.Range("L2:L6").WrapText = $true

# Row 7 now uses an explicit (smaller) custom height
$allies.Rows.Item(7).RowHeight = 12.75

$allies.Range("K18").Select()
$bands.Activate()

Write-Output "Band.xlsx data refresh applied"
